$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.433.15'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +9.01%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.602.09'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +8.29%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.71%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +9.39%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9923'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.96%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3689'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.74%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3402'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +10.48%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.72'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.75%  '

# Row 10
$ws.Range('E10').Value = '  +7.55%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07058'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.78%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9992'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.59%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.945'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +7.72%  '

# Row 14
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.68'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +8.95%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.643'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.92%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001087'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +5.54%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.600.25'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +8.33%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9923'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.83%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06666'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +12.22%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.07'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +12.41%  '

# Row 21
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +11.16%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.037'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +9.99%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.84'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +7.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.432.42'
$ws.Range('D24').ClearFormats()

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.401'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.24%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.516'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +18.00%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '150.36'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +5.81%  '

# Row 28
$ws.Range('E28').Value = '  +12.96%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.780.07'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.69%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.70'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.77%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.200'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +6.41%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.040'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +20.28%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9532'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +16.17%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08262'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.34%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.640'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +6.44%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.307'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +12.58%  '

# Row 37
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.281'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.45%  '

# Row 38
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '11.89'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +13.71%  '

# Row 39
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.595'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +12.32%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06122'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.93%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02221'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +8.67%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2032'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +7.72%  '

# Row 43
$ws.Range('E43').Value = '  +1.79%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5922'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +11.70%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.863'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +9.33%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.24'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.11%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5704'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +9.71%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.87'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.21%  '

# Row 49
$ws.Range('E49').Value = '  +9.09%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06817'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.13%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.96'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +9.44%  '

